$d = $word.ActiveDocument

function Set-ParagraphXml($para, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

# --- Edit 1: split "but it was felt..." run into "as" + " it was felt..." --------------
# (paragraph: "We have taken the definition of done ... as part of a task being done.")
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*but it was felt it was important to highlight the need for*") {
        $p1 = $p
        break
    }
}

$p1body = @'
<w:p>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">We have taken the definition of done to mean a story </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">is fully implemented into production. Within our Gitlab project we have a &#8220;Done&#8221; </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>column</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> for tasks that have met all criteria for being done except being implemented</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> into master. </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">The definition has changed little since it was first defined in week 1 </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>as</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> it was felt it was important to highlight the need for </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">a </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>written code review</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> as part of a task being done. </w:t></w:r>
</w:p>
'@
Set-ParagraphXml $p1 $p1body

# --- Edit 2: "valuable as part of the next sprint" -> "most-valuable for the next sprint",
#     split across several runs ------------------------------------------------------------
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*enter the sprint-backlog if it*deemed to be valuable as part of the next sprint*") {
        $p2 = $p
        break
    }
}

$p2body = @'
<w:p>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">We have taken the definition of ready to mean when a task is ready to </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">enter the sprint-backlog if it&#8217;s deemed to be </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>most-</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">valuable </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>for the</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> next sprint. </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">The Definition of ready has </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">seen some improvements seen </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>it&#8217;s</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> initial inception in week 1. Mainly the need for tasks to follow INVEST criteria in particular tasks being small and independent. </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>The team added clarification for outlining tasks</w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">to improve team understanding of what&#8217;s expected. </w:t></w:r>
</w:p>
'@
Set-ParagraphXml $p2 $p2body

# --- Edit 3: turn the empty "ind left=1080" ListParagraph into a numbered list item -----
$p3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "List Paragraph" -and $p.Format.LeftIndent -eq 54) {
        $p3 = $p
        break
    }
}

$p3body = @'
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">The story should be understandable </w:t></w:r>
<w:r><w:t xml:space="preserve">to </w:t></w:r>
<w:r><w:t xml:space="preserve">a client as </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>it&#8217;s</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> purpose is to invite conversation with the client. </w:t></w:r>
</w:p>
'@
Set-ParagraphXml $p3 $p3body
